$p = $ppt.ActivePresentation

# --- 1) Handout master date placeholder: "17.1.2023 г." -> "18.01.23 г." ---
$hm = $p.HandoutMaster
$hmDate = $hm.Shapes.Item(2)
$hmDate.TextFrame.TextRange.Text = "18.01.23 г."

# --- 2) Notes master date placeholder: "17-Jan-23" -> "1/18/23" ---
$nm = $p.NotesMaster
$nmDate = $nm.Shapes.Item(2)
$nmDate.TextFrame.TextRange.Text = "1/18/23"

# --- 3) Slide 11 ("Method ToString()"): resize the code textbox ---
$s11 = $p.Slides.Item(11)
$codeBox = $s11.Shapes.Item(4)
$codeBox.Top = 224.5348031496063
$codeBox.Width = 797.2270866141732

# --- 4) Slide 11: drop the trailing space in "...name is " ---
$tr = $codeBox.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf("is {FirstName}")
$spacePos = $idx + 3
$tr.Characters($spacePos, 1).Text = ""
